$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "48.986.68"
Set-TextValue $ws.Range("E2") "  +1.90%  "
Set-TextValue $ws.Range("D3") "2.531.38"
Set-TextValue $ws.Range("E3") "  +0.94%  "
Set-TextValue $ws.Range("D4") "0.998"
Set-TextValue $ws.Range("E4") "  -0.22%  "
Set-TextValue $ws.Range("D5") "320.08"
Set-TextValue $ws.Range("E5") "  -0.12%  "
Set-TextValue $ws.Range("D6") "107.13"
Set-TextValue $ws.Range("E6") "  -1.52%  "
Set-TextValue $ws.Range("D7") "0.525"
Set-TextValue $ws.Range("E7") "  -0.73%  "
Set-TextValue $ws.Range("D8") "0.998"
Set-TextValue $ws.Range("E8") "  -0.15%  "
Set-TextValue $ws.Range("D9") "0.548"
Set-TextValue $ws.Range("E9") "  +0.70%  "
Set-TextValue $ws.Range("D10") "39.63"
Set-TextValue $ws.Range("E10") "  -0.46%  "
Set-TextValue $ws.Range("D11") "20.14"
Set-TextValue $ws.Range("E11") "  +0.23%  "
Set-TextValue $ws.Range("D12") "0.0808"
Set-TextValue $ws.Range("E12") "  -1.11%  "
Set-TextValue $ws.Range("D13") "0.126"
Set-TextValue $ws.Range("E13") "  +0.89%  "
Set-TextValue $ws.Range("D14") "7.18"
Set-TextValue $ws.Range("E14") "  -0.16%  "
Set-TextValue $ws.Range("D15") "2.924.05"
Set-TextValue $ws.Range("E15") "  +0.78%  "
Set-TextValue $ws.Range("D16") "2.556.21"
Set-TextValue $ws.Range("E16") "  +0.23%  "
Set-TextValue $ws.Range("D17") "0.849"
Set-TextValue $ws.Range("E17") "  +0.42%  "
Set-TextValue $ws.Range("D18") "48.778.02"
Set-TextValue $ws.Range("E18") "  +1.78%  "
Set-TextValue $ws.Range("D19") "13.00"
Set-TextValue $ws.Range("E19") "  -1.44%  "
Set-TextValue $ws.Range("E20") "  +8.32%  "
Set-TextValue $ws.Range("D21") "6.64"
Set-TextValue $ws.Range("E21") "  +0.69%  "
Set-TextValue $ws.Range("D22") "0.0₃0939"
Set-TextValue $ws.Range("E22") "  -0.29%  "
Set-TextValue $ws.Range("D23") "282.39"
Set-TextValue $ws.Range("E23") "  +2.99%  "
Set-TextValue $ws.Range("D24") "71.23"
Set-TextValue $ws.Range("E24") "  -1.29%  "
Set-TextValue $ws.Range("E25") "  -1.91%  "
Set-TextValue $ws.Range("D26") "26.14"
Set-TextValue $ws.Range("E26") "  +1.04%  "
Set-TextValue $ws.Range("E28") "  -7.57%  "
Set-TextValue $ws.Range("D29") "0.144"
Set-TextValue $ws.Range("E29") "  +2.32%  "
Set-TextValue $ws.Range("D30") "9.73"
Set-TextValue $ws.Range("E30") "  -3.27%  "
Set-TextValue $ws.Range("D31") "35.06"
Set-TextValue $ws.Range("E31") "  -0.76%  "
Set-TextValue $ws.Range("D32") "49.57"
Set-TextValue $ws.Range("E32") "  +0.21%  "
Set-TextValue $ws.Range("D33") "19.50"
Set-TextValue $ws.Range("E33") "  +1.03%  "
Set-TextValue $ws.Range("E34") "  -0.19%  "
Set-TextValue $ws.Range("D35") "5.32"
Set-TextValue $ws.Range("E35") "  -0.33%  "
Set-TextValue $ws.Range("D36") "0.0778"
Set-TextValue $ws.Range("E36") "  -0.83%  "
Set-TextValue $ws.Range("D37") "1.99"
Set-TextValue $ws.Range("E37") "  +1.66%  "
Set-TextValue $ws.Range("D38") "4.61"
Set-TextValue $ws.Range("E38") "  -0.12%  "
Set-TextValue $ws.Range("D39") "2.94"
Set-TextValue $ws.Range("E39") "  -0.77%  "
Set-TextValue $ws.Range("E40") "  -0.42%  "
Set-TextValue $ws.Range("D41") "2.21"
Set-TextValue $ws.Range("E41") "  -0.06%  "
Set-TextValue $ws.Range("D44") "0.0305"
Set-TextValue $ws.Range("E44") "  -0.13%  "
Set-TextValue $ws.Range("E47") "  +7.00%  "
Set-TextValue $ws.Range("E48") "  +5.78%  "
Set-TextValue $ws.Range("D49") "9.01"
Set-TextValue $ws.Range("E49") "  -0.25%  "
Set-TextValue $ws.Range("D50") "5.22"
Set-TextValue $ws.Range("E50") "  +0.80%  "
Set-TextValue $ws.Range("D51") "80.77"
Set-TextValue $ws.Range("E51") "  +1.66%  "

# Row 42/43 swap: Monero <-> EnergySwap
Set-TextValue $ws.Range("B42") "EnergySwap"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D42") "22.13"
Set-TextValue $ws.Range("E42") "  +1.44%  "

Set-TextValue $ws.Range("B43") "Monero"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D43") "119.80"
Set-TextValue $ws.Range("E43") "  -2.27%  "

# Row 45/46 swap: Maker <-> NEARProtocol
Set-TextValue $ws.Range("B45") "NEARProtocol"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D45") "3.24"
Set-TextValue $ws.Range("E45") "  +3.92%  "

Set-TextValue $ws.Range("B46") "Maker"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D46") "2.006.47"
Set-TextValue $ws.Range("E46") "  -0.68%  "
